$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H18").Value = "Digikey"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4.6100000000000003
$ws.Range("E18").Value = 9.2200000000000006
$ws.Range("E23").Value = 41.02

$ws.Range("J18").Select()
